$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Invalid (G) and Absent (H) -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count (D) and Real (E) -> 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Rows 5-18: Absent (H) -> 1
for ($r = 5; $r -le 18; $r++) {
    $ws.Range("H$r").Value = 1
}
